$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.917.53"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").Value = "1.633.35"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'216.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("D6").Value = "'0.5132"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.71%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "'0.2569"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").Value = "'0.06351"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.25%  "

$ws.Range("D10").Value = "'19.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.42%  "

$ws.Range("D11").Value = "'0.07775"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("D12").Value = "'4.254"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.80%  "

$ws.Range("D13").Value = "1.636.05"
$ws.Range("E13").Value = "  -0.19%  "

$ws.Range("D14").Value = "1.858.13"
$ws.Range("E14").Value = "  -0.56%  "

$ws.Range("D15").Value = "'0.5524"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.41%  "

$ws.Range("D16").Value = "'63.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.55%  "

$ws.Range("D17").Value = "0.0₅7630"
$ws.Range("E17").Value = "  -1.51%  "

$ws.Range("D18").Value = "25.938.15"
$ws.Range("E18").Value = "  -0.45%  "

$ws.Range("D19").Value = "'1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").Value = "'195.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.82%  "

$ws.Range("D21").Value = "'4.430"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("D22").Value = "'9.866"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.67%  "

$ws.Range("D23").Value = "'6.030"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("E24").Value = "  -0.31%  "

$ws.Range("D25").Value = "'1.893"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.08%  "

$ws.Range("D26").Value = "'142.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.54%  "

$ws.Range("D27").Value = "'0.1263"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.81%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'15.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.20%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'6.764"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.02%  "

$ws.Range("E30").Value = "  +0.53%  "

$ws.Range("D31").Value = "'0.04918"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.93%  "

$ws.Range("D32").Value = "'3.233"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.84%  "

$ws.Range("D33").Value = "'3.190"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.61%  "

$ws.Range("D34").Value = "'1.547"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.40%  "

$ws.Range("D35").Value = "'2.373"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.50%  "

$ws.Range("D36").Value = "'0.8983"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("D37").Value = "'0.5527"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.11%  "

$ws.Range("D38").Value = "'2.536"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.60%  "

$ws.Range("D39").Value = "1.116.03"
$ws.Range("E39").Value = "  -2.24%  "

$ws.Range("D40").Value = "'0.01559"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.31%  "

$ws.Range("D41").Value = "'0.9999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("D42").Value = "'5.583"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.40%  "

$ws.Range("D43").Value = "'0.7963"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.64%  "

$ws.Range("D44").Value = "'97.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.45%  "

$ws.Range("D45").Value = "1.767.69"
$ws.Range("E45").Value = "  -0.65%  "

$ws.Range("D46").Value = "0.0₈118"
$ws.Range("E46").Value = "  -7.45%  "

$ws.Range("D47").Value = "'0.4431"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.05%  "

$ws.Range("D48").Value = "'1.002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.16%  "

$ws.Range("D49").Value = "'54.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("D50").Value = "'0.05133"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.49%  "

$ws.Range("D51").Value = "'7.558"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.21%  "
